$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 18:35"

# --- Refresh daily COVID counters for a handful of existing countries ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1733543
$ws.Cells.Item(4, 3).Value = 8268
$ws.Cells.Item(4, 4).Value = 481367
$ws.Cells.Item(4, 5).Value = 1151137
$ws.Cells.Item(4, 7).Value = 467
$ws.Cells.Item(4, 8).Value = 101039

# Row 9: Italia
$ws.Cells.Item(9, 2).Value = 231139
$ws.Cells.Item(9, 3).Value = 584
$ws.Cells.Item(9, 4).Value = 147101
$ws.Cells.Item(9, 5).Value = 50966
$ws.Cells.Item(9, 7).Value = 117
$ws.Cells.Item(9, 8).Value = 33072

# Row 13: India
$ws.Cells.Item(13, 2).Value = 157484
$ws.Cells.Item(13, 3).Value = 6691
$ws.Cells.Item(13, 4).Value = 67285
$ws.Cells.Item(13, 5).Value = 85676
$ws.Cells.Item(13, 7).Value = 179
$ws.Cells.Item(13, 8).Value = 4523

# Row 56: Argelia
$ws.Cells.Item(56, 2).Value = 8857
$ws.Cells.Item(56, 3).Value = 160
$ws.Cells.Item(56, 4).Value = 5129
$ws.Cells.Item(56, 5).Value = 3105
$ws.Cells.Item(56, 7).Value = 6
$ws.Cells.Item(56, 8).Value = 623

# --- Re-sort a cluster of low-count countries; their country label moves to a
# different row, together with the "Casos activos" (D) / "Muertes" (H) figures
# that belong to that country ---

# Rows 199-201: Belice / Nueva Caledonia / Santa Lucia trio re-ordered
$ws.Cells.Item(199, 1).Value = "Santa Lucia"
$ws.Cells.Item(199, 4).Value = 18
$ws.Cells.Item(199, 8).Value = 0

$ws.Cells.Item(200, 1).Value = "Belice"
$ws.Cells.Item(200, 4).Value = 16
$ws.Cells.Item(200, 8).Value = 2

$ws.Cells.Item(201, 1).Value = "Nueva Caledonia"

# Rows 207-208: Islas Turcas y Caicos / Groenlandia swap
$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(207, 4).Value = 11
$ws.Cells.Item(207, 8).Value = 0

$ws.Cells.Item(208, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208, 4).Value = 10
$ws.Cells.Item(208, 8).Value = 1

# Rows 210-211: Seychelles / Montserrat swap
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0
